# "Generate Report for Handback" — fills in the Latest Target File / Latest
# Handback File / Latest Handback DateTime / Error Detail columns for the
# e4b2ccb0-ae52-4dc1-ac84-3730fe297376 row on both the zh-cn and de-de
# report sheets, widens the J/K/R columns to match the other "wide" columns,
# and hyperlinks the new Latest Target File cell back to the source .md file.

$wb = $excel.ActiveWorkbook

# Excel's ColumnWidth (characters) is stored in the OOXML as
# characters + 5/6 (the standard cell-padding constant), so subtract that
# off to land on an exact raw width of 40.
$targetColWidth = 40 - (5/6)

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4e7a0c1c1ddf44fb22b59417275f1ae3afa4b557/e2e/e4b2ccb0-ae52-4dc1-ac84-3730fe297376.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3c700bd5eb7bf2266a9b3d93b36a56b789de324c/e2e/e4b2ccb0-ae52-4dc1-ac84-3730fe297376.md."
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3c700bd5eb7bf2266a9b3d93b36a56b789de324c/e2e/e4b2ccb0-ae52-4dc1-ac84-3730fe297376.md"

$sheetsInfo = @(
    @{ Name = "zh-cn"; Handback = "e4b2ccb0-ae52-4dc1-ac84-3730fe297376.2e03ce94d2a0f8bfcafebaf48063cf1b757c05b6.zh-cn.xlf"; DateTime = "2017-02-17 07:39:46" },
    @{ Name = "de-de"; Handback = "e4b2ccb0-ae52-4dc1-ac84-3730fe297376.2e03ce94d2a0f8bfcafebaf48063cf1b757c05b6.de-de.xlf"; DateTime = "2017-02-17 07:40:12" }
)

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Widen Latest Target File (J), Latest Handback File (K) and Error Detail (R)
    $ws.Columns.Item(10).ColumnWidth = $targetColWidth
    $ws.Columns.Item(11).ColumnWidth = $targetColWidth
    $ws.Columns.Item(18).ColumnWidth = $targetColWidth

    # Latest Target File (J6) -> hyperlink to the source markdown file
    $ws.Range("J6").Value = "e4b2ccb0-ae52-4dc1-ac84-3730fe297376.md"
    $ws.Range("J6").Font.Underline = 2
    $ws.Range("J6").Font.Color = 15570276
    $ws.Hyperlinks.Add($ws.Range("J6"), $targetUrl, "", "", "e4b2ccb0-ae52-4dc1-ac84-3730fe297376.md") | Out-Null

    # Latest Handback File (K6)
    $ws.Range("K6").Value = $info.Handback

    # Latest Handback DateTime (L6)
    $ws.Range("L6").Value = $info.DateTime

    # Error Detail (R6)
    $ws.Range("R6").Value = $errorDetail
}
